# Backlog.xlsx edit: mark several tasks as "terminado" and refresh the
# autofilter so it also shows "en proceso" items; the filtered range grows
# from A1:C100 to A1:C104, and the last-used selection moves to B76.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark tasks as completed ("terminado")
$ws.Range("B40").Value = "terminado"
$ws.Range("B96").Value = "terminado"
$ws.Range("B99").Value = "terminado"
$ws.Range("B100").Value = "terminado"
$ws.Range("B102").Value = "terminado"
$ws.Range("B103").Value = "terminado"

# Refresh the autofilter over the full data range, now also showing
# "en proceso" tasks in addition to "no comenzado"
$ws.AutoFilterMode = $false
$ws.Range("A1:C104").AutoFilter(2, @("en proceso", "no comenzado"), 7)

# Row 40 stays visible even though its status became "terminado"
$ws.Rows.Item(40).Hidden = $false

# Update the hidden _FilterDatabase defined name to the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$C`$104"
    }
}

# Move the active selection ahead of printing
$ws.Activate()
$ws.Range("B76").Select()
